$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for rows 2-51.
# NumberFormat is forced to text ("@") before assigning D-column values
# because several price strings (e.g. "1.000", "273.40", "0.2908") are
# valid numeric literals and Excel/COM would otherwise silently coerce
# them to Double and normalize away meaningful trailing zeros / dot
# formatting used by this sheet (European-style grouped numbers).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.631.99"
$ws.Range("E2").Value = "  +1.05%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.873.61"
$ws.Range("E3").Value = "  +0.20%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.51"
$ws.Range("E5").Value = "  +0.98%  "

$ws.Range("E6").Value = "  -0.01%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4728"
$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2908"
$ws.Range("E8").Value = "  +1.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06484"
$ws.Range("E9").Value = "  +0.22%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.12"
$ws.Range("E10").Value = "  +4.74%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07706"
$ws.Range("E11").Value = "  -0.91%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7376"
$ws.Range("E12").Value = "  +1.13%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.46"
$ws.Range("E13").Value = "  +1.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.871.59"
$ws.Range("E14").Value = "  -0.02%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.161"
$ws.Range("E15").Value = "  +0.25%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "273.40"
$ws.Range("E16").Value = "  -0.61%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.652.73"
$ws.Range("E17").Value = "  +1.09%  "

$ws.Range("E18").Value = "  -0.41%  "

$ws.Range("E19").Value = "  +0.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007520"
$ws.Range("E20").Value = "  -0.48%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.117.56"
$ws.Range("E21").Value = "  +0.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9989"
$ws.Range("E22").Value = "  -0.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.264"
$ws.Range("E23").Value = "  +0.31%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.193"
$ws.Range("E24").Value = "  +0.30%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.214"
$ws.Range("E25").Value = "  -0.55%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.84"
$ws.Range("E26").Value = "  -1.16%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.75"
$ws.Range("E27").Value = "  -0.84%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.909"
$ws.Range("E28").Value = "  -0.39%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1002"
$ws.Range("E29").Value = "  +1.24%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.343"
$ws.Range("E30").Value = "  -2.70%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.509"
$ws.Range("E31").Value = "  -0.96%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.280"
$ws.Range("E32").Value = "  -0.71%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.101"
$ws.Range("E33").Value = "  +1.34%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04799"
$ws.Range("E34").Value = "  +0.38%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.120"
$ws.Range("E35").Value = "  -0.16%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6962"
$ws.Range("E36").Value = "  -0.33%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.718"
$ws.Range("E37").Value = "  -0.13%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01853"
$ws.Range("E38").Value = "  +0.38%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.758"
$ws.Range("E39").Value = "  +0.42%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.244"
$ws.Range("E40").Value = "  -2.52%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.17"
$ws.Range("E41").Value = "  +4.12%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.972"
$ws.Range("E42").Value = "  +2.83%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4177"
$ws.Range("E43").Value = "  +1.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8358"
$ws.Range("E45").Value = "  -0.68%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.22"
$ws.Range("E46").Value = "  -0.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.332"
$ws.Range("E47").Value = "  -0.23%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.44"
$ws.Range("E48").Value = "  +0.59%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.989"
$ws.Range("E49").Value = "  -1.60%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "918.37"
$ws.Range("E50").Value = "  +0.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05651"
$ws.Range("E51").Value = "  +1.40%  "
